$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update individual data cells per the diff
$ws.Range("D2").Value = 2
$ws.Range("G2").Value = 1
$ws.Range("G4").Value = 4
$ws.Range("I7").Value = 8
$ws.Range("I8").Value = 4
$ws.Range("D12").Value = 1
$ws.Range("H12").Value = 12

# I14 becomes a formula instead of a hardcoded value
$ws.Range("I14").Formula = "=SUM(I2:I13)"

# Update the sheet view: remove topLeftCell override and change selection
$ws.Activate()
$ws.Range("N19").Select()
